$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("mitprofessors", "October", "Addition", "Zufferey, Raphael", "Assistant Professor (Starting January 2025)"),
    @("michiganprofessors", "October", "Addition", "Jacinto Ulloa", "Assistant Professor, Mechanical Engineering"),
    @("michiganprofessors", "October", "Deletion", "Jesse Austin-Breneman", "Assistant Professor, Mechanical Engineering"),
    @("UIUCprofessors", "October", "Addition", "Anthony Jacobi", "Department Head, Richard W. Kritzer Distinguished Professor"),
    @("UIUCprofessors", "October", "Addition", "Cunjiang Yu", "Founder Professor"),
    @("UIUCprofessors", "October", "Deletion", "Tony Jacobi", "Department Head, Richard W. Kritzer Distinguished Professor"),
    @("georgiaprofessors", "October", "Addition", "Christopher J. Saldaña", "Ring Family Professor"),
    @("georgiaprofessors", "October", "Addition", "Samuel Graham", "Professor")
)

$startRow = 54
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
